$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for rows 2-29 from 45448 to 45449
$ws.Range("C2:C29").Value = 45449

# Row 29 specific updates: new case number, updated date, updated area
$ws.Range("A29").Value = "A 22929-2024"
$ws.Range("B29").Value = 45448
$ws.Range("G29").Value = 7.3
